$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range("E2").Value = "2026-02-04 19:21:04"
$ws.Range("H2").Value = "'84%"
$ws.Range("E3").Value = "2026-02-04 19:21:07"
$ws.Range("E4").Value = "2026-02-04 19:21:09"
$ws.Range("O4").Value = "6.2 °C"
$ws.Range("E5").Value = "2026-02-04 19:21:12"
$ws.Range("J5").Value = "992.6 hPa"
$ws.Range("O5").Value = "8.7 °C"
$ws.Range("E6").Value = "2026-02-04 19:21:15"
$ws.Range("E7").Value = "2026-02-04 19:21:17"
$ws.Range("H7").Value = "'82%"
$ws.Range("E8").Value = "2026-02-04 19:21:20"
$ws.Range("H8").Value = "'86%"
$ws.Range("O8").Value = "7.6 °C"
$ws.Range("E9").Value = "2026-02-04 19:21:23"
$ws.Range("E10").Value = "2026-02-04 19:21:25"
$ws.Range("E11").Value = "2026-02-04 19:21:28"
$ws.Range("E12").Value = "2026-02-04 19:21:31"
$ws.Range("O12").Value = "8.5 °C"
$ws.Range("E13").Value = "2026-02-04 19:21:34"
$ws.Range("E14").Value = "2026-02-04 19:21:36"
$ws.Range("O14").Value = "-6.1 °C"
$ws.Range("E15").Value = "2026-02-04 19:21:39"
$ws.Range("E16").Value = "2026-02-04 19:21:42"
$ws.Range("O16").Value = "2.9 °C"
$ws.Range("E17").Value = "2026-02-04 19:21:44"
$ws.Range("E18").Value = "2026-02-04 19:21:47"
$ws.Range("O18").Value = "-6.8 °C"
$ws.Range("E19").Value = "2026-02-04 19:21:50"
$ws.Range("J19").Value = "994.3 hPa"
$ws.Range("E20").Value = "2026-02-04 19:21:53"
$ws.Range("O20").Value = "-4.5 °C"
$ws.Range("E21").Value = "2026-02-04 19:21:56"
$ws.Range("H21").Value = "'75%"
$ws.Range("O21").Value = "6.0 °C"
$ws.Range("E22").Value = "2026-02-04 19:21:59"
$ws.Range("E23").Value = "2026-02-04 19:22:01"
$ws.Range("H23").Value = "'73%"
$ws.Range("O23").Value = "8.7 °C"
$ws.Range("E24").Value = "2026-02-04 19:22:04"
$ws.Range("E25").Value = "2026-02-04 19:22:07"
$ws.Range("O25").Value = "1.3 °C"
$ws.Range("E26").Value = "2026-02-04 19:22:10"
$ws.Range("H26").Value = "'67%"
$ws.Range("O26").Value = "-1.8 °C"
$ws.Range("E27").Value = "2026-02-04 19:22:12"
$ws.Range("E28").Value = "2026-02-04 19:22:15"
$ws.Range("O28").Value = "2.6 °C"
$ws.Range("E29").Value = "2026-02-04 19:22:18"
$ws.Range("E30").Value = "2026-02-04 19:22:21"
$ws.Range("H30").Value = "'75%"
$ws.Range("E31").Value = "2026-02-04 19:22:23"
$ws.Range("E32").Value = "2026-02-04 19:22:25"
$ws.Range("E33").Value = "2026-02-04 19:22:28"
$ws.Range("O33").Value = "9.6 °C"
$ws.Range("E34").Value = "2026-02-04 19:22:31"
$ws.Range("E35").Value = "2026-02-04 19:22:34"
$ws.Range("E36").Value = "2026-02-04 19:22:36"
